$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5891093015670776
$ws.Range("B1").Value = 0.9424312114715576
$ws.Range("C1").Value = 2.715994596481323
$ws.Range("D1").Value = 6.13934850692749
$ws.Range("E1").Value = 2.090463399887085
